$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Clientes": trim the registros table down to a single data row and
# fix up a couple of values (this is the "login" data the commit references).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Clientes")

# Drop the duplicate/garbage rows 3-6, keeping only the header and row 2.
$ws.Rows("3:6").Delete()

# Replace row 2's contents with the corrected record. The ID and phone
# number look numeric but must be stored as text, so format those cells as
# Text before assigning the values.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "1034918141"
$ws.Cells.Item(2, 2).Value = "Daniel"
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "50516800"

# Narrow the Telefono column slightly (12 -> 10 characters).
$ws.Columns.Item(3).ColumnWidth = 9.14

# ---------------------------------------------------------------------------
# Sheet "Usuarios": the login fix itself - the stored password was the
# numeric placeholder 123; replace it with the real text password "hola".
# ---------------------------------------------------------------------------
$wsUsuarios = $wb.Worksheets.Item("Usuarios")
$wsUsuarios.Cells.Item(2, 3).Value = "hola"
